$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 75

# A75 holds a date-looking label ("01-04-2021") that must stay a plain text
# string (matching the existing column of text dates) rather than being
# auto-converted into a date serial number by the smart-text parser.
# Enter it as a formula returning the literal string, then convert the
# formula to its literal value in place (Copy + PasteSpecial values) so the
# stored cell is a plain shared string with no special number formatting.
$ws.Range("A75").Formula = '="01-04-2021"'
$ws.Range("A75").Copy()
$ws.Range("A75").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Cells.Item($row, 2).Value = 1265
$ws.Cells.Item($row, 3).Value = 0
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 1265
$ws.Cells.Item($row, 7).Value = 12
$ws.Cells.Item($row, 8).Value = 88
$ws.Cells.Item($row, 9).Value = 1165
